# Add 8 new rows (55-62) describing Huawei Atlas AI accelerator cards
# (Atlas 300I Model 3010, Atlas 300I Model 9000, Atlas 300I Pro, Atlas 300V Pro)
# for both aarch64 and x86_64 builds of openEuler 20.03 LTS, plus their
# download/support hyperlinks in column Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row=55; A="19e5"; B="d100"; C="0200"; D="0100"; E="aarch64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I(Model 3010)"; O="Atlas 300I(Model 3010)"; P=$null; Q="https://support.huawei.com/enterprise/zh/ascend-computing/a300-3010-pid-251560253/software" },
    @{ Row=56; A="19e5"; B="d100"; C="0200"; D="0100"; E="x86_64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I(Model 3010)"; O="Atlas 300I(Model 3010)"; P=$null; Q="https://support.huawei.com/enterprise/zh/ascend-computing/a300-3010-pid-251560253/software" },
    @{ Row=57; A="19e5"; B="d801"; C="0200"; D="0100"; E="aarch64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I(Model 9000)"; O="Atlas 300I(Model 9000)"; P=$null; Q="https://support.huawei.com/enterprise/zh/ascend-computing/a300t-9000-pid-250702906/software" },
    @{ Row=58; A="19e5"; B="d801"; C="0200"; D="0100"; E="x86_64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I(Model 9000)"; O="Atlas 300I(Model 9000)"; P=$null; Q="https://support.huawei.com/enterprise/zh/ascend-computing/a300t-9000-pid-250702906/software" },
    @{ Row=59; A="19e5"; B="d500"; C="0200"; D="0100"; E="aarch64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I Pro"; O="Atlas 300I Pro"; P="02313FUJ"; Q="https://support.huawei.com/enterprise/zh/ascend-computing/atlas-300i-pro-pid-251052354/software" },
    @{ Row=60; A="19e5"; B="d500"; C="0200"; D="0100"; E="x86_64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300I Pro"; O="Atlas 300I Pro"; P="02313FUJ"; Q="https://support.huawei.com/enterprise/zh/ascend-computing/atlas-300i-pro-pid-251052354/software" },
    @{ Row=61; A="19e5"; B="d500"; C="0200"; D="0100"; E="aarch64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300V Pro"; O="Atlas 300V Pro"; P="02313SHF/02314BCJ"; Q="https://support.huawei.com/enterprise/zh/ascend-computing/atlas-300i-pro-pid-251052354/software" },
    @{ Row=62; A="19e5"; B="d500"; C="0200"; D="0100"; E="x86_64"; F="openEuler 20.03 LTS"; G="npu"; I="AI"; J="2022.09.21"; M="Huawei"; N="Atlas 300V Pro"; O="Atlas 300V Pro"; P="02313SHF/02314BCJ"; Q="https://support.huawei.com/enterprise/zh/ascend-computing/atlas-300i-pro-pid-251052354/software" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    [void]($ws.Cells.Item($rowNum, 1).Value = $r.A)   # vendorID
    [void]($ws.Cells.Item($rowNum, 2).Value = $r.B)   # deviceID
    [void]($ws.Cells.Item($rowNum, 3).Value = $r.C)   # svID
    [void]($ws.Cells.Item($rowNum, 4).Value = $r.D)   # ssID
    [void]($ws.Cells.Item($rowNum, 5).Value = $r.E)   # architecture
    [void]($ws.Cells.Item($rowNum, 6).Value = $r.F)   # os
    [void]($ws.Cells.Item($rowNum, 7).Value = $r.G)   # driverName
    [void]($ws.Cells.Item($rowNum, 9).Value = $r.I)   # type
    [void]($ws.Cells.Item($rowNum, 10).Value = $r.J)  # date
    [void]($ws.Cells.Item($rowNum, 13).Value = $r.M)  # chipVendor
    [void]($ws.Cells.Item($rowNum, 14).Value = $r.N)  # boardModel
    [void]($ws.Cells.Item($rowNum, 15).Value = $r.O)  # chipModel

    if ($r.P) {
        [void]($ws.Cells.Item($rowNum, 16).Value = $r.P)  # item
    }

    $qCell = $ws.Cells.Item($rowNum, 17)
    [void]($qCell.Value = $r.Q)
    [void]$ws.Hyperlinks.Add($qCell, $r.Q)
}
